$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# 1) Split the trailing "...come into issues with their work." run so the
#    final period becomes its own run (mirrors the diff's run split).
$len = $tr.Length
$lastChar = $tr.Characters($len, 1)
$lastChar.Text = "."

# 2) Append a new bullet paragraph with the additional reason text directly
#    after the "work." paragraph.
$apostrophe = [string]([char]0x2019)
$newText = "Doesn" + $apostrophe + "t rely on synchronization of all clients."
$tr.InsertAfter("`r" + $newText)

# 3) Insert a blank bullet paragraph between the "work." paragraph and the
#    new reason paragraph (same bullet/indent formatting is inherited).
$tr2 = $tf.TextRange
$newPara = $tr2.Paragraphs(6, 1)
$newPara.InsertBefore("`r")
